$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws1 = $wb.Worksheets.Item("Forecast Comparison")
Set-TextValue $ws1.Range("B2") "2025-01-12"
$ws1.Range("D2").Value = 2
Set-TextValue $ws1.Range("B3") "2025-01-19"
$ws1.Range("D3").Value = 2
Set-TextValue $ws1.Range("B4") "2025-01-26"
$ws1.Range("D4").Value = 2
Set-TextValue $ws1.Range("B5") "2025-02-02"
$ws1.Range("D5").Value = 2
Set-TextValue $ws1.Range("B6") "2025-02-09"
$ws1.Range("D6").Value = 2
Set-TextValue $ws1.Range("B7") "2025-02-16"
$ws1.Range("D7").Value = 2
Set-TextValue $ws1.Range("B8") "2025-02-23"
$ws1.Range("D8").Value = 3
Set-TextValue $ws1.Range("B9") "2025-03-02"
$ws1.Range("D9").Value = 3
Set-TextValue $ws1.Range("B10") "2025-03-09"
$ws1.Range("D10").Value = 3
Set-TextValue $ws1.Range("B11") "2025-03-16"
$ws1.Range("D11").Value = 3
Set-TextValue $ws1.Range("B12") "2025-03-23"
$ws1.Range("D12").Value = 3
Set-TextValue $ws1.Range("B13") "2025-03-30"
$ws1.Range("D13").Value = 3
Set-TextValue $ws1.Range("B14") "2025-04-06"
$ws1.Range("D14").Value = 3
Set-TextValue $ws1.Range("B15") "2025-04-13"
$ws1.Range("D15").Value = 2
Set-TextValue $ws1.Range("B16") "2025-04-20"
$ws1.Range("D16").Value = 2
Set-TextValue $ws1.Range("B17") "2025-04-27"
$ws1.Range("D17").Value = 2

$ws2 = $wb.Worksheets.Item("Summary")
Set-TextValue $ws2.Range("B2") "2022-12-25 to 2025-01-05"
Set-TextValue $ws2.Range("B8") "657 units"
Set-TextValue $ws2.Range("B9") "41"
Set-TextValue $ws2.Range("B10") "20"
Set-TextValue $ws2.Range("B11") "9"
Set-TextValue $ws2.Range("B12") "3"
Set-TextValue $ws2.Range("B13") "2025-03-02"
Set-TextValue $ws2.Range("B14") "2"
Set-TextValue $ws2.Range("B15") "2025-04-27"
